$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments ---
$ws.Columns.Item(3).ColumnWidth = 24.17   # C: 24 -> 25
$ws.Columns.Item(6).ColumnWidth = 19.17   # F: 22 -> 20
$ws.Columns.Item(11).ColumnWidth = 24.17  # K: 24 -> 25

# --- Cell value / formatting updates (per this week's ride assignments) ---
$ws.Range("E2").Value = 'Driver: Olivia Chang'
$ws.Range("F2").Value = 'Driver: Grace Park'
$ws.Range("C3").Value = 'Lindsey Ro'
$ws.Range("D3").Value = 'Joann Jung'
$ws.Range("E3").Value = 'Kyle Hwang'
$ws.Range("F3").Value = 'Chae Moon'
$ws.Range("M3").Value = 'Ariel Bahn'
$ws.Range("O3").Value = 'Zoe Li'
$ws.Range("C4").Value = 'Ariel Bahn'
$ws.Range("D4").Value = 'Isabelle Li'
$ws.Range("E4").Value = 'Sam Ko'
$ws.Range("F4").Value = 'Samuel Wen'
$ws.Range("G4").Value = 'Jane Yoo'
$ws.Range("K4").Value = 'Jane Yoo'
$ws.Range("L4").Value = 'Daniel Kim'
$ws.Range("M4").Value = 'Jasmine Cheng'
$ws.Range("O4").Value = 'Claire Doh'
$ws.Range("C5").Value = 'Zoe Li'
$ws.Range("D5").Value = 'Phillip Seo'
$ws.Range("E5").Value = 'Aaron duong'
$ws.Range("F5").Value = 'Gabriel Ni'
$ws.Range("G5").Value = 'Grace Sowon Park'
$ws.Range("K5").Value = 'Grace Sowon Park'
$ws.Range("L5").Value = 'Daniel Kuo'
$ws.Range("N5").Value = 'Isabelle Li'
$ws.Range("O5").Value = 'Sam Ko'
$ws.Range("D6").Value = 'Stella Son'
$ws.Range("F6").Value = 'Cara Lee'
$ws.Range("N6").Value = 'Phillip Seo'
$ws.Range("O6").Value = 'Cara Lee'
$ws.Range("C9").Value = 'Driver: Kaitlyn Kim'
$ws.Range("D9").Value = 'Driver: Josh Paik'
$ws.Range("K9").Value = 'Driver: Claudia Parents'
$ws.Range("C10").Value = 'Karina Pan'
$ws.Range("C10").Interior.Color = 15597567
$ws.Range("D10").Value = 'Ben Kim'
$ws.Range("D10").Interior.Color = 13421812
$ws.Range("E10").Value = 'Claire Doh'
$ws.Range("C11").Value = 'Faith Chen'
$ws.Range("C11").Interior.Color = 15597567
$ws.Range("D11").Value = 'Jay Jung'
$ws.Range("D11").Interior.Color = 13421812
$ws.Range("E11").Value = 'Joel Shim'
$ws.Range("M11").Value = 'Irene Jung'
$ws.Range("N11").Value = 'Lindsey Ro'
$ws.Range("C12").Clear()
$ws.Range("D12").Value = 'Daniel Kim'
$ws.Range("D12").Interior.Color = 13421812
$ws.Range("M12").Value = 'Stella Son'
$ws.Range("N12").Value = 'Christina Ko'
$ws.Range("C13").Clear()
$ws.Range("D13").Value = 'Daniel Kuo'
$ws.Range("D13").Interior.Color = 13421812
$ws.Range("E13").Value = 'Irene Jung'
$ws.Range("M13").Value = 'Kyle Hwang'
$ws.Range("C16").Value = 'Driver: Claudia Parents'
$ws.Range("E17").Value = 'Christina Ko — No valid driver'
$ws.Range("M17").Value = 'Samuel Wen — No valid driver'
$ws.Range("E18").Value = 'Jasmine Cheng — No valid driver'
$ws.Range("E19").Value = 'Susanna Tang — No valid driver'
$ws.Range("M19").Value = 'Eugene Seo — No valid driver'
$ws.Range("E20").Value = 'Eugene Seo — No valid driver'
$ws.Range("M20").Value = 'helena song🐟 — No valid driver'
$ws.Range("E21").Value = 'Josh Yang — No valid driver'
$ws.Range("E22").Value = 'helena song🐟 — No valid driver'
$ws.Range("M22").Value = 'Khang Le — No valid driver'
$ws.Range("E23").Value = 'Darius Ajebon — No valid driver'
$ws.Range("M23").Value = 'Jacob Lei — No valid driver'
$ws.Range("E24").Value = 'Jacob Lei — No valid driver'
